# Update points 79174419 -> 0.00
# Adds a new trailing row (row 14) to Sheet1 for phone "79174419" with a
# blank birthday and 0 total_points, matching the existing rows' layout
# (phone stored as text, birthday blank-text, total_points numeric 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

# Leading "'" forces text entry (matches the sheet's text-typed phone/
# birthday columns) without leaving the cell's stored text empty.
$ws.Range("A$row").Value = "'79174419"
$ws.Range("B$row").Value = "'"
$ws.Range("C$row").Value = 0

# Drop the "quote prefix" formatting that typing a leading apostrophe
# applies, so the new cells carry no extra style (same as the rest of
# the sheet's data rows).
$ws.Range("A$row`:B$row").ClearFormats()
